$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "35.208.52"  # D2: '35.207.71' -> '35.208.52'
$ws.Cells.Item(2, 5).Value = "  +0.21%  "  # E2: '  +0.17%  ' -> '  +0.21%  '
$ws.Cells.Item(3, 4).Value = "1.889.96"  # D3: '1.890.49' -> '1.889.96'
$ws.Cells.Item(3, 5).Value = "  +1.98%  "  # E3: '  +1.99%  ' -> '  +1.98%  '
$ws.Cells.Item(4, 5).Value = "  -0.28%  "  # E4: '  -0.23%  ' -> '  -0.28%  '
$ws.Cells.Item(5, 4).Value = "'242.71"  # D5: '242.57' -> '242.71'
$ws.Cells.Item(5, 5).Value = "  +1.99%  "  # E5: '  +1.98%  ' -> '  +1.99%  '
$ws.Cells.Item(6, 4).Value = "'0.652"  # D6: '0.653' -> '0.652'
$ws.Cells.Item(6, 5).Value = "  +5.09%  "  # E6: '  +5.15%  ' -> '  +5.09%  '
$ws.Cells.Item(7, 5).Value = "  -0.28%  "  # E7: '  -0.27%  ' -> '  -0.28%  '
$ws.Cells.Item(8, 4).Value = "'41.13"  # D8: '41.09' -> '41.13'
$ws.Cells.Item(8, 5).Value = "  -1.69%  "  # E8: '  -1.80%  ' -> '  -1.69%  '
$ws.Cells.Item(9, 5).Value = "  +5.88%  "  # E9: '  +5.87%  ' -> '  +5.88%  '
$ws.Cells.Item(10, 4).Value = "'50.00"  # D10: '50.07' -> '50.00'
$ws.Cells.Item(10, 5).Value = "  +7.49%  "  # E10: '  +7.63%  ' -> '  +7.49%  '
$ws.Cells.Item(11, 5).Value = "  +2.03%  "  # E11: '  +1.97%  ' -> '  +2.03%  '
$ws.Cells.Item(12, 5).Value = "  +0.66%  "  # E12: '  +0.71%  ' -> '  +0.66%  '
$ws.Cells.Item(13, 4).Value = "2.164.68"  # D13: '2.167.27' -> '2.164.68'
$ws.Cells.Item(13, 5).Value = "  +1.99%  "  # E13: '  +2.17%  ' -> '  +1.99%  '
$ws.Cells.Item(14, 4).Value = "'11.86"  # D14: '11.87' -> '11.86'
$ws.Cells.Item(14, 5).Value = "  +4.21%  "  # E14: '  +4.30%  ' -> '  +4.21%  '
$ws.Cells.Item(15, 4).Value = "1.905.32"  # D15: '1.897.91' -> '1.905.32'
$ws.Cells.Item(15, 5).Value = "  +2.17%  "  # E15: '  +3.31%  ' -> '  +2.17%  '
$ws.Cells.Item(16, 5).Value = "  +2.27%  "  # E16: '  +2.26%  ' -> '  +2.27%  '
$ws.Cells.Item(17, 5).Value = "  +2.40%  "  # E17: '  +2.20%  ' -> '  +2.40%  '
$ws.Cells.Item(18, 4).Value = "35.190.21"  # D18: '35.213.28' -> '35.190.21'
$ws.Cells.Item(18, 5).Value = "  +0.19%  "  # E18: '  +0.21%  ' -> '  +0.19%  '
$ws.Cells.Item(19, 4).Value = "'71.11"  # D19: '71.05' -> '71.11'
$ws.Cells.Item(19, 5).Value = "  +1.65%  "  # E19: '  +1.59%  ' -> '  +1.65%  '
$ws.Cells.Item(21, 4).Value = "'240.18"  # D21: '240.31' -> '240.18'
$ws.Cells.Item(21, 5).Value = "  -0.15%  "  # E21: '  -0.09%  ' -> '  -0.15%  '
$ws.Cells.Item(22, 4).Value = "'12.37"  # D22: '12.36' -> '12.37'
$ws.Cells.Item(22, 5).Value = "  +1.47%  "  # E22: '  +1.44%  ' -> '  +1.47%  '
$ws.Cells.Item(23, 4).Value = "'4.73"  # D23: '4.72' -> '4.73'
$ws.Cells.Item(23, 5).Value = "  +0.32%  "  # E23: '  +0.15%  ' -> '  +0.32%  '
$ws.Cells.Item(24, 5).Value = "  -0.28%  "  # E24: '  -0.30%  ' -> '  -0.28%  '
$ws.Cells.Item(25, 4).Value = "'2.41"  # D25: '2.43' -> '2.41'
$ws.Cells.Item(25, 5).Value = "  +33.02%  "  # E25: '  +32.82%  ' -> '  +33.02%  '
$ws.Cells.Item(26, 4).Value = "'2.29"  # D26: '2.28' -> '2.29'
$ws.Cells.Item(26, 5).Value = "  +0.66%  "  # E26: '  +0.37%  ' -> '  +0.66%  '
$ws.Cells.Item(27, 4).Value = "'169.79"  # D27: '169.91' -> '169.79'
$ws.Cells.Item(27, 5).Value = "  +0.44%  "  # E27: '  +0.48%  ' -> '  +0.44%  '
$ws.Cells.Item(28, 5).Value = "  +5.26%  "  # E28: '  +5.21%  ' -> '  +5.26%  '
$ws.Cells.Item(29, 4).Value = "'18.15"  # D29: '18.16' -> '18.15'
$ws.Cells.Item(29, 5).Value = "  +3.11%  "  # E29: '  +3.26%  ' -> '  +3.11%  '
$ws.Cells.Item(30, 5).Value = "  +2.08%  "  # E30: '  +2.24%  ' -> '  +2.08%  '
$ws.Cells.Item(31, 5).Value = "  +2.93%  "  # E31: '  +2.85%  ' -> '  +2.93%  '
$ws.Cells.Item(32, 2).Value = "ImmutableX"  # B32: 'BinanceUSD' -> 'ImmutableX'
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"  # C32: 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' -> 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).Value = "'0.939"  # D32: '1.01' -> '0.939'
$ws.Cells.Item(32, 5).Value = "  +14.35%  "  # E32: '  -0.20%  ' -> '  +14.35%  '
$ws.Cells.Item(33, 2).Value = "BinanceUSD"  # B33: 'Hedera' -> 'BinanceUSD'
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"  # C33: 'https://coinranking.com/coin/jad286TjB+hedera-hbar' -> 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(33, 4).Value = "'1.01"  # D33: '0.0558' -> '1.01'
$ws.Cells.Item(33, 5).Value = "  -0.25%  "  # E33: '  +0.67%  ' -> '  -0.25%  '
$ws.Cells.Item(34, 2).Value = "Hedera"  # B34: 'ImmutableX' -> 'Hedera'
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"  # C34: 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' -> 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = "'0.0558"  # D34: '0.935' -> '0.0558'
$ws.Cells.Item(34, 5).Value = "  +0.64%  "  # E34: '  +15.84%  ' -> '  +0.64%  '
$ws.Cells.Item(35, 4).Value = "'4.09"  # D35: '4.08' -> '4.09'
$ws.Cells.Item(35, 5).Value = "  +2.39%  "  # E35: '  +1.85%  ' -> '  +2.39%  '
$ws.Cells.Item(36, 5).Value = "  -1.81%  "  # E36: '  -1.69%  ' -> '  -1.81%  '
$ws.Cells.Item(37, 5).Value = "  +0.29%  "  # E37: '  +0.42%  ' -> '  +0.29%  '
$ws.Cells.Item(38, 5).Value = "  +0.40%  "  # E38: '  +0.88%  ' -> '  +0.40%  '
$ws.Cells.Item(39, 4).Value = "'0.0208"  # D39: '0.0207' -> '0.0208'
$ws.Cells.Item(39, 5).Value = "  +3.74%  "  # E39: '  +3.54%  ' -> '  +3.74%  '
$ws.Cells.Item(40, 4).Value = "'1.09"  # D40: '1.08' -> '1.09'
$ws.Cells.Item(40, 5).Value = "  +1.47%  "  # E40: '  +1.43%  ' -> '  +1.47%  '
$ws.Cells.Item(41, 4).Value = "'0.0640"  # D41: '0.0637' -> '0.0640'
$ws.Cells.Item(41, 5).Value = "  +16.03%  "  # E41: '  +15.15%  ' -> '  +16.03%  '
$ws.Cells.Item(42, 4).Value = "'15.93"  # D42: '15.90' -> '15.93'
$ws.Cells.Item(42, 5).Value = "  +7.21%  "  # E42: '  +7.09%  ' -> '  +7.21%  '
$ws.Cells.Item(43, 4).Value = "'88.74"  # D43: '88.75' -> '88.74'
$ws.Cells.Item(43, 5).Value = "  -1.11%  "  # E43: '  -1.13%  ' -> '  -1.11%  '
$ws.Cells.Item(44, 4).Value = "1.334.08"  # D44: '1.335.02' -> '1.334.08'
$ws.Cells.Item(44, 5).Value = "  -0.39%  "  # E44: '  -0.42%  ' -> '  -0.39%  '
$ws.Cells.Item(45, 2).Value = "MultiversX"  # B45: 'RenderToken' -> 'MultiversX'
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"  # C45: 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' -> 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(45, 4).Value = "'47.96"  # D45: '2.35' -> '47.96'
$ws.Cells.Item(45, 5).Value = "  +38.77%  "  # E45: '  +1.82%  ' -> '  +38.77%  '
$ws.Cells.Item(46, 2).Value = "RenderToken"  # B46: 'MultiversX' -> 'RenderToken'
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"  # C46: 'https://coinranking.com/coin/omwkOTglq+multiversx-egld' -> 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(46, 4).Value = "'2.34"  # D46: '47.81' -> '2.34'
$ws.Cells.Item(46, 5).Value = "  +1.58%  "  # E46: '  +38.62%  ' -> '  +1.58%  '
$ws.Cells.Item(47, 5).Value = "  -1.33%  "  # E47: '  -1.26%  ' -> '  -1.33%  '
$ws.Cells.Item(48, 4).Value = "'2.77"  # D48: '2.76' -> '2.77'
$ws.Cells.Item(48, 5).Value = "  +1.17%  "  # E48: '  +0.77%  ' -> '  +1.17%  '
$ws.Cells.Item(49, 5).Value = "  +0.61%  "  # E49: '  +0.83%  ' -> '  +0.61%  '
$ws.Cells.Item(50, 4).Value = "2.075.05"  # D50: '2.076.83' -> '2.075.05'
$ws.Cells.Item(50, 5).Value = "  +1.64%  "  # E50: '  +1.83%  ' -> '  +1.64%  '
$ws.Cells.Item(51, 4).Value = "'11.21"  # D51: '11.16' -> '11.21'
$ws.Cells.Item(51, 5).Value = "  -12.90%  "  # E51: '  -14.00%  ' -> '  -12.90%  '
